$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (they use "."
# both as decimal point and as a thousands separator, e.g.
# "25.836.86" or "215.54"). Plain assignment would let Excel
# auto-convert single-dot values into floating point numbers, so we
# briefly force a text number format while writing the value, then
# restore the default "Normal" style so the cell keeps matching the
# surrounding (unstyled) data cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.836.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.629.14'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.59%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5104'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06405'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  -1.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.267'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.28%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.629.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.853.83'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5602'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.48'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7556'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.837.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.62%  '

$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.334'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.63%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.810'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.007'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.25%  '

$ws.Range("E24").Value = '  +0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.823'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.75%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1287'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '141.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.753'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.43%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.79%  '

$ws.Range("E30").Value = '  -0.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04891'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.312'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.191'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.559'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.379'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8964'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.132.85'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.545'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5493'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01560'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.004'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.32%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.588'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.26%  '

$ws.Range("E43").Value = '  -1.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.36'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.776.57'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.10%  '

$ws.Range("E46").Value = '  -7.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4427'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05066'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.562'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.000'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.60%  '
